$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New order of model names for rows 2..26 (column A)
$names = @(
    "model_32_6_0",
    "model_32_6_22",
    "model_32_6_21",
    "model_32_6_20",
    "model_32_6_19",
    "model_32_6_18",
    "model_32_6_17",
    "model_32_6_16",
    "model_32_6_15",
    "model_32_6_14",
    "model_32_6_13",
    "model_32_6_23",
    "model_32_6_12",
    "model_32_6_10",
    "model_32_6_9",
    "model_32_6_8",
    "model_32_6_7",
    "model_32_6_6",
    "model_32_6_5",
    "model_32_6_4",
    "model_32_6_3",
    "model_32_6_2",
    "model_32_6_1",
    "model_32_6_11",
    "model_32_6_24"
)

for ($i = 0; $i -lt $names.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 1).Value = $names[$i]
}

# New shared metric values (identical across all 25 data rows, columns B..Q)
$valB = [double]"0.9999949039827802"
$valC = [double]"0.9991177146462858"
$valD = [double]"0.9999976561244449"
$valE = [double]"0.9999999738209685"
$valF = [double]"0.9999992291165175"
$valG = [double]"4.756908385316126e-06"
$valH = [double]"0.0008235746498410597"
$valI = [double]"1.292685094922034e-06"
$valJ = [double]"9.113317272946016e-09"
$valK = [double]"6.508992355609909e-07"
$valL = [double]"8.583594873526996e-05"
$valM = [double]"0.002181033788210565"
$valN = [double]"1.000004892176531"
$valO = [double]"0.00227388491193704"
$valP = [double]"122.5118251992398"
$valQ = [double]"182.2367406177816"

for ($r = 2; $r -le 26; $r++) {
    $ws.Cells.Item($r, 2).Value = $valB
    $ws.Cells.Item($r, 3).Value = $valC
    $ws.Cells.Item($r, 4).Value = $valD
    $ws.Cells.Item($r, 5).Value = $valE
    $ws.Cells.Item($r, 6).Value = $valF
    $ws.Cells.Item($r, 7).Value = $valG
    $ws.Cells.Item($r, 8).Value = $valH
    $ws.Cells.Item($r, 9).Value = $valI
    $ws.Cells.Item($r, 10).Value = $valJ
    $ws.Cells.Item($r, 11).Value = $valK
    $ws.Cells.Item($r, 12).Value = $valL
    $ws.Cells.Item($r, 13).Value = $valM
    $ws.Cells.Item($r, 14).Value = $valN
    $ws.Cells.Item($r, 15).Value = $valO
    $ws.Cells.Item($r, 16).Value = $valP
    $ws.Cells.Item($r, 17).Value = $valQ
}
